# Updates cryptos list values (Price / Volume(1h) columns) per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.646.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.286.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.48%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.503"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.27%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.500"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0788"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "50.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.77%  "
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.65"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.642.39"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.275.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.787"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.566.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0895"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.99%  "
$ws.Range("E22").Value = "  -3.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "234.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.22%  "
$ws.Range("E25").Value = "  -3.03%  "
$ws.Range("E26").Value = "  -3.12%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "165.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.37%  "
$ws.Range("E30").Value = "  -8.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.61%  "
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("E34").Value = "  -3.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.39"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0695"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.49%  "
$ws.Range("E37").Value = "  -5.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.22"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.91%  "
$ws.Range("E40").Value = "  -6.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.57%  "
$ws.Range("E42").Value = "  -2.49%  "
$ws.Range("E43").Value = "  -1.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.962.19"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0282"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.84%  "
$ws.Range("E47").Value = "  -7.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.87%  "
$ws.Range("E50").Value = "  -4.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.507.86"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.10%  "
